# Update the Dsc2-Dsg1a LR-pairs sheet with newly recomputed TPM values.
# The new data only has 3 result rows (ECs/FAPs/MuSCs all paired against
# MuSCs as target cluster) instead of the previous 6 rows, so the extra
# rows are removed and rows 2-4 get the refreshed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-obsolete rows 5-7 first so the used range shrinks to A1:T4.
$ws.Rows("5:7").Delete() | Out-Null

# Row 2: ECs -> Dsc2 -> Dsg1a -> MuSCs
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Dsc2"
$ws.Range("C2").Value2 = "Dsg1a"
$ws.Range("D2").Value2 = "MuSCs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.09662633333333333
$ws.Range("H2").Value2 = 0.289879
$ws.Range("I2").Value2 = 0.1167569426878839
$ws.Range("J2").Value2 = 0.1167569426878839
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.011155
$ws.Range("N2").Value2 = 0.033465
$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 0.001077866748333333
$ws.Range("R2").Value2 = 0.009700800735000001
$ws.Range("S2").Value2 = 0.1167569426878839
$ws.Range("T2").Value2 = 0.1167569426878839

# Row 3: FAPs -> Dsc2 -> Dsg1a -> MuSCs
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Dsc2"
$ws.Range("C3").Value2 = "Dsg1a"
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.01499333333333333
$ws.Range("H3").Value2 = 0.04498
$ws.Range("I3").Value2 = 0.01811696356790599
$ws.Range("J3").Value2 = 0.01811696356790599
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.011155
$ws.Range("N3").Value2 = 0.033465
$ws.Range("O3").Value2 = 1
$ws.Range("P3").Value2 = 1
$ws.Range("Q3").Value2 = 0.0001672506333333333
$ws.Range("R3").Value2 = 0.0015052557
$ws.Range("S3").Value2 = 0.01811696356790599
$ws.Range("T3").Value2 = 0.01811696356790599

# Row 4: MuSCs -> Dsc2 -> Dsg1a -> MuSCs
$ws.Range("A4").Value2 = "MuSCs"
$ws.Range("B4").Value2 = "Dsc2"
$ws.Range("C4").Value2 = "Dsg1a"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.7159656666666666
$ws.Range("H4").Value2 = 2.147897
$ws.Range("I4").Value2 = 0.86512609374421
$ws.Range("J4").Value2 = 0.86512609374421
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.011155
$ws.Range("N4").Value2 = 0.033465
$ws.Range("O4").Value2 = 1
$ws.Range("P4").Value2 = 1
$ws.Range("Q4").Value2 = 0.007986597011666666
$ws.Range("R4").Value2 = 0.07187937310500001
$ws.Range("S4").Value2 = 0.86512609374421
$ws.Range("T4").Value2 = 0.86512609374421
